# Update cryptos list cell values per the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.955.21"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "3.408.97"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.66%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.408.93"
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.122"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").Value = "3.995.80"
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "3.412.54"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("E16").Value = "  -4.01%  "
$ws.Range("D17").Value = "63.005.54"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.560"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("D24").Value = "3.548.09"
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("E27").Value = "  -7.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "  -5.24%  "
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("E31").Value = "  -4.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.153"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("D34").Value = "3.438.22"
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "163.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0764"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.782"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.61%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("E46").Value = "  -5.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("E50").Value = "  -3.72%  "
$ws.Range("D51").Value = "2.263.73"
$ws.Range("E51").Value = "  -5.30%  "
